$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the merged tuple-like string combining the former A2:A4 values
$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"

# Remove the now-obsolete rows 3 and 4 (previously held the type line and P/T)
$ws.Rows("3:4").Delete()
